$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 21940
$ws.Range("J3").Value = 21940
$ws.Range("L3").Value = 21940
$ws.Range("N3").Value = -22168
$ws.Range("H17").Value = 1210.2273
$ws.Range("I17").Value = 890
$ws.Range("J17").Value = 1225.4762
$ws.Range("K17").Value = 2670
$ws.Range("L17").Value = 3676.4286
$ws.Range("M17").Value = -2502
$ws.Range("N17").Value = -4012.4286
$ws.Range("H33").Value = 692.6326
$ws.Range("I33").Value = 716.76746
$ws.Range("J33").Value = 519.6667
$ws.Range("K33").Value = 716.76746
$ws.Range("L33").Value = 519.6667
$ws.Range("M33").Value = -487.76746
$ws.Range("N33").Value = -977.6667
$ws.Range("H43").Value = 1279.625
$ws.Range("I43").Value = 1196.8889
$ws.Range("J43").Value = 1386
$ws.Range("K43").Value = 1196.8889
$ws.Range("L43").Value = 1386
$ws.Range("M43").Value = -1127.8889
$ws.Range("N43").Value = -1524
$ws.Range("H94").Value = 0
$ws.Range("I94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("M94").ClearContents()
$ws.Range("H102").Value = 21940
$ws.Range("J102").Value = 21940
$ws.Range("L102").Value = 21940
$ws.Range("N102").Value = -28430
$ws.Range("H134").Value = 200038000
$ws.Range("J134").Value = 200038000
$ws.Range("L134").Value = 200038000
$ws.Range("N134").Value = -200048140
$ws.Range("H138").Value = 4325.3184
$ws.Range("I138").Value = 3062.4666
$ws.Range("J138").Value = 4978.517
$ws.Range("K138").Value = 9187.399800000001
$ws.Range("L138").Value = 14935.551
$ws.Range("M138").Value = -4047.399800000001
$ws.Range("N138").Value = -25215.551

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1695.9048
$ws.Range("I2").Value = 1251.6154
$ws.Range("J2").Value = 2417.875
$ws.Range("K2").Value = 1251.6154
$ws.Range("L2").Value = 2417.875
$ws.Range("M2").Value = -1138.6154
$ws.Range("N2").Value = -2643.875
$ws.Range("H32").Value = 2116580
$ws.Range("I32").Value = 15610.728
$ws.Range("J32").Value = 71448570
$ws.Range("K32").Value = 15610.728
$ws.Range("L32").Value = 71448570
$ws.Range("M32").Value = -15323.728
$ws.Range("N32").Value = -71449144
$ws.Range("H45").Value = 920.25
$ws.Range("I45").Value = 744.6923
$ws.Range("J45").Value = 1072.4
$ws.Range("K45").Value = 744.6923
$ws.Range("L45").Value = 1072.4
$ws.Range("M45").Value = -367.6923
$ws.Range("N45").Value = -1826.4
$ws.Range("H74").Value = 1649.1765
$ws.Range("I74").Value = 1516.8572
$ws.Range("J74").Value = 2266.6667
$ws.Range("K74").Value = 1516.8572
$ws.Range("L74").Value = 2266.6667
$ws.Range("M74").Value = -642.8571999999999
$ws.Range("N74").Value = -4014.6667
$ws.Range("H77").Value = 1649.1765
$ws.Range("I77").Value = 1516.8572
$ws.Range("J77").Value = 2266.6667
$ws.Range("K77").Value = 7584.286
$ws.Range("L77").Value = 11333.3335
$ws.Range("M77").Value = -3216.286
$ws.Range("N77").Value = -20069.3335
$ws.Range("H110").Value = 911.37933
$ws.Range("I110").Value = 516.64703
$ws.Range("J110").Value = 1470.5834
$ws.Range("K110").Value = 516.64703
$ws.Range("L110").Value = 1470.5834
$ws.Range("M110").Value = 1528.35297
$ws.Range("N110").Value = -5560.5834
$ws.Range("H116").Value = 1695.9048
$ws.Range("I116").Value = 1251.6154
$ws.Range("J116").Value = 2417.875
$ws.Range("K116").Value = 1251.6154
$ws.Range("L116").Value = 2417.875
$ws.Range("M116").Value = 1042.3846
$ws.Range("N116").Value = -7005.875
$ws.Range("H122").Value = 1138.8096
$ws.Range("I122").Value = 847.4
$ws.Range("J122").Value = 1403.7273
$ws.Range("K122").Value = 2542.2
$ws.Range("L122").Value = 4211.1819
$ws.Range("M122").Value = -92.19999999999982
$ws.Range("N122").Value = -9111.1819
$ws.Range("H132").Value = 9008.172
$ws.Range("I132").Value = 12836.692
$ws.Range("J132").Value = 6745.864
$ws.Range("K132").Value = 38510.076
$ws.Range("L132").Value = 20237.592
$ws.Range("M132").Value = -35980.076
$ws.Range("N132").Value = -25297.592

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1695.9048
$ws.Range("I3").Value = 1251.6154
$ws.Range("J3").Value = 2417.875
$ws.Range("K3").Value = 1251.6154
$ws.Range("L3").Value = 2417.875
$ws.Range("M3").Value = -1137.6154
$ws.Range("N3").Value = -2645.875
$ws.Range("H26").Value = 22083.334
$ws.Range("I26").Value = 14300
$ws.Range("J26").Value = 61000
$ws.Range("K26").Value = 14300
$ws.Range("L26").Value = 61000
$ws.Range("M26").Value = -14008
$ws.Range("N26").Value = -61584
$ws.Range("H97").Value = 8324.182000000001
$ws.Range("I97").Value = 1230.8334
$ws.Range("J97").Value = 16836.2
$ws.Range("K97").Value = 1230.8334
$ws.Range("L97").Value = 16836.2
$ws.Range("M97").Value = -239.8334
$ws.Range("N97").Value = -18818.2
$ws.Range("H105").Value = 2583.2104
$ws.Range("I105").Value = 2451.8235
$ws.Range("J105").Value = 3700
$ws.Range("K105").Value = 2451.8235
$ws.Range("L105").Value = 3700
$ws.Range("M105").Value = -704.8235
$ws.Range("N105").Value = -7194

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 83334700
$ws.Range("I122").Value = 100000850
$ws.Range("K122").Value = 300002550
$ws.Range("M122").Value = -300000100

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H59").Value = 1006.25
$ws.Range("H70").Value = 5150
$ws.Range("J70").Value = 5934.706
$ws.Range("L70").Value = 17804.118
$ws.Range("N70").Value = -18434.118
$ws.Range("H73").Value = 5150
$ws.Range("J73").Value = 5934.706
$ws.Range("L73").Value = 17804.118
$ws.Range("N73").Value = -19988.118
$ws.Range("H137").Value = 4540.2
$ws.Range("I137").Value = 3307.5
$ws.Range("J137").Value = 4699.2583
$ws.Range("K137").Value = 9922.5
$ws.Range("L137").Value = 14097.7749
$ws.Range("M137").Value = -4822.5
$ws.Range("N137").Value = -24297.7749

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 8902.526
$ws.Range("I122").Value = 15183.667
$ws.Range("J122").Value = 3249.5
$ws.Range("K122").Value = 45551.001
$ws.Range("L122").Value = 9748.5
$ws.Range("M122").Value = -43101.001
$ws.Range("N122").Value = -14648.5
$ws.Range("H132").Value = 3217.2122
$ws.Range("I132").Value = 2518.8
$ws.Range("J132").Value = 3799.2222
$ws.Range("K132").Value = 7556.400000000001
$ws.Range("L132").Value = 11397.6666
$ws.Range("M132").Value = -5026.400000000001
$ws.Range("N132").Value = -16457.6666

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1492.4667
$ws.Range("I122").Value = 1651.7142
$ws.Range("J122").Value = 1120.8889
$ws.Range("K122").Value = 4955.142599999999
$ws.Range("L122").Value = 3362.6667
$ws.Range("M122").Value = -2505.142599999999
$ws.Range("N122").Value = -8262.6667
